$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("REMOTE PM")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Backup"

$new.Range("H3:M6").UnMerge()
$new.Range("H8:M10").UnMerge()
$new.Range("H12:M16").UnMerge()

$new.Rows("27:37").Delete()

$new.Range("H3:M5").Merge()
$new.Range("H7:M9").Merge()
